$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1818287.8
$ws.Range("I5").Value = 2000108.4
$ws.Range("J5").Value = 81
$ws.Range("K5").Value = 2000108.4
$ws.Range("L5").Value = 81
$ws.Range("M5").Value = -1999993.4
$ws.Range("N5").Value = -311
$ws.Range("H8").Value = 35883
$ws.Range("I8").Value = 41830.25
$ws.Range("J8").Value = 199.5
$ws.Range("K8").Value = 125490.75
$ws.Range("L8").Value = 598.5
$ws.Range("M8").Value = -125351.75
$ws.Range("N8").Value = -876.5
$ws.Range("H17").Value = 1775.3684
$ws.Range("J17").Value = 1775.3684
$ws.Range("L17").Value = 5326.1052
$ws.Range("N17").Value = -5662.1052
$ws.Range("H95").Value = 62996.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 62996.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 62996.5
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -68488.5
$ws.Range("H111").Value = 331
$ws.Range("I111").Value = 341.33334
$ws.Range("K111").Value = 1024.00002
$ws.Range("M111").Value = 2042.99998
$ws.Range("H129").Value = 3418.7778
$ws.Range("I129").Value = 3281.2856
$ws.Range("K129").Value = 9843.856800000001
$ws.Range("M129").Value = -4843.856800000001
$ws.Range("H132").Value = 371708.12
$ws.Range("I132").Value = 1404.84
$ws.Range("J132").Value = 5000499
$ws.Range("K132").Value = 4214.52
$ws.Range("L132").Value = 15001497
$ws.Range("M132").Value = -1684.52
$ws.Range("N132").Value = -15006557
$ws.Range("H137").Value = 3148.8708
$ws.Range("I137").Value = 2804.8125
$ws.Range("K137").Value = 8414.4375
$ws.Range("M137").Value = -5864.4375
$ws.Range("H138").Value = 2896.8635
$ws.Range("J138").Value = 2510.3381
$ws.Range("L138").Value = 7531.0143
$ws.Range("N138").Value = -17811.0143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H63").Value = 4166.3335
$ws.Range("I63").Value = 1999.5
$ws.Range("J63").Value = 5249.75
$ws.Range("K63").Value = 1999.5
$ws.Range("L63").Value = 5249.75
$ws.Range("M63").Value = -1313.5
$ws.Range("N63").Value = -6621.75
$ws.Range("H66").Value = 4166.3335
$ws.Range("I66").Value = 1999.5
$ws.Range("J66").Value = 5249.75
$ws.Range("K66").Value = 9997.5
$ws.Range("L66").Value = 26248.75
$ws.Range("M66").Value = -6565.5
$ws.Range("N66").Value = -33112.75
$ws.Range("H122").Value = 2613
$ws.Range("I122").Value = 2613
$ws.Range("K122").Value = 7839
$ws.Range("M122").Value = -5389

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 44893
$ws.Range("J81").Value = 44893
$ws.Range("L81").Value = 44893
$ws.Range("N81").Value = -47015
$ws.Range("H84").Value = 44893
$ws.Range("J84").Value = 44893
$ws.Range("L84").Value = 134679
$ws.Range("N84").Value = -145287
$ws.Range("H107").Value = 1252.8889
$ws.Range("I107").Value = 1170.4286
$ws.Range("J107").Value = 1541.5
$ws.Range("K107").Value = 1170.4286
$ws.Range("L107").Value = 1541.5
$ws.Range("M107").Value = 749.5714
$ws.Range("N107").Value = -5381.5
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 22934.285
$ws.Range("J28").Value = 22934.285
$ws.Range("L28").Value = 22934.285
$ws.Range("N28").Value = -23424.285
$ws.Range("H41").Value = 11384.538
$ws.Range("J41").Value = 11818.091
$ws.Range("L41").Value = 11818.091
$ws.Range("N41").Value = -12674.091
$ws.Range("H74").Value = 69749
$ws.Range("J74").Value = 69749
$ws.Range("L74").Value = 69749
$ws.Range("N74").Value = -71497
$ws.Range("H77").Value = 69749
$ws.Range("J77").Value = 69749
$ws.Range("L77").Value = 209247
$ws.Range("N77").Value = -217983
$ws.Range("H93").Value = 19543.555
$ws.Range("I93").Value = 8699
$ws.Range("K93").Value = 8699
$ws.Range("M93").Value = -6827
$ws.Range("H105").Value = 1796.875
$ws.Range("I105").Value = 1645.1
$ws.Range("K105").Value = 1645.1
$ws.Range("M105").Value = 101.9000000000001
$ws.Range("H122").Value = 1279.75
$ws.Range("I122").Value = 1374.6666
$ws.Range("K122").Value = 4123.9998
$ws.Range("M122").Value = -1673.9998
$ws.Range("H134").Value = 2943.4348
$ws.Range("I134").Value = 2866.8333
$ws.Range("K134").Value = 8600.499899999999
$ws.Range("M134").Value = -6065.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 377.69232
$ws.Range("I12").Value = 473
$ws.Range("K12").Value = 1419
$ws.Range("M12").Value = -1246
$ws.Range("H47").Value = 471.2
$ws.Range("J47").Value = 862.5
$ws.Range("L47").Value = 2587.5
$ws.Range("N47").Value = -3449.5
$ws.Range("H92").Value = 250.93333
$ws.Range("I92").Value = 186.25
$ws.Range("K92").Value = 558.75
$ws.Range("M92").Value = 689.25
$ws.Range("H107").Value = 1801.9
$ws.Range("J107").Value = 1777.375
$ws.Range("L107").Value = 5332.125
$ws.Range("N107").Value = -9172.125
$ws.Range("H109").Value = 2399.0833
$ws.Range("I109").Value = 1923.4
$ws.Range("K109").Value = 5770.200000000001
$ws.Range("M109").Value = -4730.200000000001
$ws.Range("H131").Value = 31932.945
$ws.Range("J131").Value = 1997.3235
$ws.Range("L131").Value = 5991.970499999999
$ws.Range("N131").Value = -16071.9705

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H70").Value = 9668.643
$ws.Range("I70").Value = 9090.143
$ws.Range("J70").Value = 10247.143
$ws.Range("K70").Value = 9090.143
$ws.Range("L70").Value = 10247.143
$ws.Range("M70").Value = -8820.143
$ws.Range("N70").Value = -10787.143
$ws.Range("H73").Value = 9668.643
$ws.Range("I73").Value = 9090.143
$ws.Range("J73").Value = 10247.143
$ws.Range("K73").Value = 9090.143
$ws.Range("L73").Value = 10247.143
$ws.Range("M73").Value = -8154.143
$ws.Range("N73").Value = -12119.143
$ws.Range("H80").Value = 8430.76
$ws.Range("I80").Value = 4758.8335
$ws.Range("J80").Value = 11820.23
$ws.Range("K80").Value = 4758.8335
$ws.Range("L80").Value = 11820.23
$ws.Range("M80").Value = -3760.8335
$ws.Range("N80").Value = -13816.23
$ws.Range("H83").Value = 8430.76
$ws.Range("I83").Value = 4758.8335
$ws.Range("J83").Value = 11820.23
$ws.Range("K83").Value = 23794.1675
$ws.Range("L83").Value = 59101.14999999999
$ws.Range("M83").Value = -18802.1675
$ws.Range("N83").Value = -69085.14999999999
$ws.Range("H102").Value = 2004.8182
$ws.Range("I102").Value = 1405.6842
$ws.Range("K102").Value = 1405.6842
$ws.Range("M102").Value = 216.3158000000001
$ws.Range("H113").Value = 4324.0835
$ws.Range("I113").Value = 3878.4
$ws.Range("K113").Value = 3878.4
$ws.Range("M113").Value = -1708.4
$ws.Range("H122").Value = 645.5
$ws.Range("I122").Value = 645.5
$ws.Range("J122").Value = 645.5
$ws.Range("K122").Value = 1936.5
$ws.Range("L122").Value = 1936.5
$ws.Range("M122").Value = 513.5
$ws.Range("N122").Value = -6836.5
$ws.Range("H126").Value = 3787.7144
$ws.Range("I126").Value = 3720.111
$ws.Range("J126").Value = 3909.4
$ws.Range("K126").Value = 11160.333
$ws.Range("L126").Value = 11728.2
$ws.Range("M126").Value = -8690.332999999999
$ws.Range("N126").Value = -16668.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3000
$ws.Range("I4").Value = 3000
$ws.Range("K4").Value = 3000
$ws.Range("M4").Value = -2887
$ws.Range("H28").Value = 3000
$ws.Range("I28").Value = 3000
$ws.Range("K28").Value = 3000
$ws.Range("M28").Value = -2768
$ws.Range("H37").Value = 3000
$ws.Range("I37").Value = 3000
$ws.Range("K37").Value = 3000
$ws.Range("M37").Value = -2893
$ws.Range("H40").Value = 4183
$ws.Range("I40").Value = 3896.7
$ws.Range("J40").Value = 4898.75
$ws.Range("K40").Value = 3896.7
$ws.Range("L40").Value = 4898.75
$ws.Range("M40").Value = -3760.7
$ws.Range("N40").Value = -5170.75
$ws.Range("H132").Value = 2618.3225
$ws.Range("I132").Value = 2470.7307
$ws.Range("J132").Value = 3385.8
$ws.Range("K132").Value = 7412.1921
$ws.Range("L132").Value = 10157.4
$ws.Range("M132").Value = -4882.1921
$ws.Range("N132").Value = -15217.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 15000
$ws.Range("J104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -21988
$ws.Range("H122").Value = 7095.2354
$ws.Range("I122").Value = 8570.25
$ws.Range("J122").Value = 3555.2
$ws.Range("K122").Value = 25710.75
$ws.Range("L122").Value = 10665.6
$ws.Range("M122").Value = -23260.75
$ws.Range("N122").Value = -15565.6
$ws.Range("H132").Value = 1157.28
$ws.Range("I132").Value = 1116.6364
$ws.Range("K132").Value = 3349.9092
$ws.Range("M132").Value = -819.9092000000001
$ws.Range("H135").Value = 66049.336
$ws.Range("J135").Value = 66049.336
$ws.Range("L135").Value = 66049.336
$ws.Range("N135").Value = -76189.336
